# Update generated figures (想去人数 / 最低票价) on the "展览" and
# "全部类型" sheets, matching the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览": F = 想去人数 (want-to-go count), G = 最低票价 (lowest price)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4454
$ws1.Range("F3").Value = 847
$ws1.Range("G4").Value = 40
$ws1.Range("G5").Value = 65
$ws1.Range("G6").Value = "不可售"
$ws1.Range("F10").Value = 617
$ws1.Range("F12").Value = 190
$ws1.Range("F13").Value = 1228
$ws1.Range("F15").Value = 2854
$ws1.Range("F16").Value = 441
$ws1.Range("F17").Value = 542

# Sheet "全部类型": same columns, mirrored data (row numbering shifts by
# one after row 10 because this sheet omits one row present in "展览").
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4454
$ws4.Range("F3").Value = 847
$ws4.Range("G4").Value = 40
$ws4.Range("G5").Value = 65
$ws4.Range("G6").Value = "不可售"
$ws4.Range("F10").Value = 617
$ws4.Range("F13").Value = 190
$ws4.Range("F14").Value = 1228
$ws4.Range("F16").Value = 2854
$ws4.Range("F17").Value = 441
$ws4.Range("F18").Value = 542
